$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '59.092.71'
$ws.Range("E2").Value = '  +2.00%  '
# Row 3
$ws.Range("D3").Value = '2.590.37'
$ws.Range("E3").Value = '  +0.85%  '
# Row 4
$ws.Range("E4").Value = '  +0.03%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '528.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.59%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.76'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.12%  '
# Row 7
$ws.Range("E7").Value = '  +0.11%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.565'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.96%  '
# Row 9
$ws.Range("D9").Value = '2.602.47'
$ws.Range("E9").Value = '  +0.70%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.44'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.13%  '
# Row 11
$ws.Range("E11").Value = '  +2.16%  '
# Row 12
$ws.Range("E12").Value = '  +1.54%  '
# Row 13
$ws.Range("E13").Value = '  +3.11%  '
# Row 14
$ws.Range("D14").Value = '3.053.35'
$ws.Range("E14").Value = '  +1.05%  '
# Row 15
$ws.Range("D15").Value = '59.031.16'
$ws.Range("E15").Value = '  +2.00%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.43'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.57%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000133'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.96%  '
# Row 18
$ws.Range("D18").Value = '2.571.96'
$ws.Range("E18").Value = '  -0.42%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '346.96'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.86%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.33'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.19%  '
# Row 21
$ws.Range("E21").Value = '  -0.01%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.39'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.75%  '
# Row 23
$ws.Range("E23").Value = '  -0.05%  '
# Row 24
$ws.Range("E24").Value = '  +3.01%  '
# Row 25
$ws.Range("E25").Value = '  +0.17%  '
# Row 26
$ws.Range("E26").Value = '  +1.74%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.16%  '
# Row 28
$ws.Range("E28").Value = '  +2.88%  '
# Row 29
$ws.Range("E29").Value = '  +0.08%  '
# Row 30
$ws.Range("E30").Value = '  +0.55%  '
# Row 31
$ws.Range("E31").Value = '  +3.03%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.85'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.61%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.77'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.86%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '148.54'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.18%  '
# Row 35
$ws.Range("E35").Value = '  +1.16%  '
# Row 36
$ws.Range("E36").Value = '  +0.02%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '36.81'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.48%  '
# Row 38
$ws.Range("E38").Value = '  +1.53%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.826'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.76%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.827'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.38%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.52'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.28%  '
# Row 42
$ws.Range("E42").Value = '  +0.07%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '268.99'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.02%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.75'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.55%  '
# Row 45
$ws.Range("E45").Value = '  +0.18%  '
# Row 46
$ws.Range("E46").Value = '  +1.84%  '
# Row 47
$ws.Range("E47").Value = '  +0.38%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.43'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.30%  '
# Row 49
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.61'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.81%  '
# Row 50
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '1.950.61'
$ws.Range("E50").Value = '  -0.42%  '
# Row 51
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0221'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.33%  '
